$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63 is currently the last data row and uses the "last row" date-only
# number format. Since we are appending a new row (64), row 63 becomes a
# regular data row and should switch to the regular datetime format used
# by the rest of the data rows (match the cell directly above it, A62).
$ws.Range("A63").NumberFormat = $ws.Range("A62").NumberFormat

# Append the new daily update row (64), with A64 taking on the special
# "last row" date-only number format that A63 used to have.
$ws.Range("A64").Value = 45804
$ws.Range("A64").NumberFormat = "YYYY-MM-DD"

$ws.Range("B64").Value = 269
$ws.Range("C64").Value = 277
$ws.Range("D64").Value = 274
